# Binomial Distribution - Cumulative Probability
# Update the inputs (number of trials n, and x) and the active selection,
# and darken the accent-coloured font used by the result cell (E8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# n (cell B2): 25 -> 20
$ws.Range("B2").Value = 20

# x / a (cell B5): 15 -> 12
$ws.Range("B5").Value = 12

# Darken the font colour used for the result (E8): Accent1 (theme) shaded
# ~-50%, i.e. RGB(31,78,121) = #1F4E79 (equivalent to theme="4"
# tint="-0.499984740745262").
$ws.Range("E8").Font.Color = 31 + 78 * 256 + 121 * 65536

# Update the active selection to D12
$ws.Range("D12").Select()
